$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.305.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.664.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5349"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.561"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.669.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.892.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8213"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.696"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.041"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1230"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.484"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05828"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.617"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.281"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.826"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5813"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8695"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.051.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.804.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4387"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.038"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.412"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.81%  "
